$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new header cells (F1, G1), copying the style of the existing header cells ---
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E1").Copy($ws.Range("G1"))

# --- The old column E held "fantasy points" data; shift that data out to the new last
#     column (G), and fill columns E (height) and F (weight) with the new scraped data ---
for ($r = 2; $r -le 16; $r++) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 243
}

# --- Set the header text for the three rightmost header cells ---
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"
